$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 100 -> 0M
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

# Row 2: 0 -> 0M
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"

# Row 3: 11 -> 0M
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# Insert 10 new rows right after row 3 (before what is currently row 4)
$anchorRow = $t.Rows.Item(4)
for ($i = 0; $i -lt 10; $i++) {
    $t.Rows.Add($anchorRow) | Out-Null
}

$newValues = @("13", "0.00003", "0.00004", "0.00003", "0.00000", "0.00003", "0.00003", "0.00004", "0.00044", "100.0")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $t.Rows.Item(4 + $i).Cells.Item(1).Range.Text = $newValues[$i]
}

# The row that used to be row 34 is now row 44 (shifted by the 10 inserted rows).
# Collapse its multi-run, tab-separated content down to a single "100".
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"

# The row that used to be row 35 is now row 45; collapse it down to "0".
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0"

# The row that used to be row 36 (empty run) is now row 46; give it the text "11".
$t.Rows.Item(46).Cells.Item(1).Range.Text = "11"
